$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Row 15 / D15: drop the stray "pen klar" fragment
# -----------------------------------------------------------------
$ws.Range("D15").ClearContents()

# -----------------------------------------------------------------
# 2. Insert 3 new rows before row 21, pushing the old
#    AIC / Metod / Inledning / Resultat / Diskussion block down
#    (old "AIC" row ends up at row 24, Metod..Diskussion end up at 25..28)
# -----------------------------------------------------------------
$ws.Rows("21:23").Insert()

# -----------------------------------------------------------------
# 3. Fill in the new / moved rows, matching the order the author
#    would naturally have typed them in (keeps new shared-string
#    ids in the same sequence as the real edit).
# -----------------------------------------------------------------

# Row 23 (new): AIC - skapa en fil med alla variabler för alla lyor
$ws.Range("A23").Value = "AIC  - skapa en fil med alla variabler för alla lyor"

# Row 24: was the original "AIC" row (now shifted down by the insert) -
# overwrite its text in place with the riplyor variant.
$ws.Range("A24").Value = "AIC  - skapa en fil med alla variabler för riplyor"

# Row 21 (new): Lyornas storlek / funkar ej / comment
$ws.Range("A21").Value = "Lyornas storlek"
$ws.Range("B21").Value = "funkar ej"
$ws.Range("C21").Value = "För svårt att se lyorna på IR-foto. Funkade inte ens att få hyfsat rätt på lyorna jag har mätt på riktigt. Skicka filen till Karin och Rasmus"

# Row 14 / C14: append extra clause to the comment text
$ws.Range("C14").Value = "Gjort shapefil och excel fil i EPSG:3006 sweref med skjutna rödrävar mellan 2000 och 2016 från Peters och Lars filer. Finns för de två senaste vintrarna på rovbasen men kommer inte åt dem. limma in i excel därifrån. det finns rödrävsreproduktion i Peters fil fram till 2008 (röd text) dock är det bara två totalt. Använd Rasmus skript. "

# Row 22 (new): Relativt mått på kullar / påbörjat / comment
$ws.Range("A22").Value = "Relativt mått på kullar"
$ws.Range("C22").Value = "Har räknat ut totala antalet kullar per lya. Har räknat ut antal inventeringar per lya mellan 2000 och 2010 (från BEBODDA_LYOR_HEF 00-10 ) och mellan 2015 och 2018 från en lista som jag plockade ut från rovbasen. Saknar alltså inventeringsdata på lyor mellan 2011 och 2014."
$ws.Range("B22").Value = "påbörjat"

$ws.Range("B23").Value = "ej påbörjat"
$ws.Range("B24").Value = "ej påbörjat"

# -----------------------------------------------------------------
# 4. Formatting: reuse the correct existing style templates so the
#    colour/weight of the "status" cells matches the other rows of
#    the same kind.
# -----------------------------------------------------------------

# B21 "funkar ej" -> bold black, same as the "avvaktar" rows (e.g. B16)
$ws.Range("B16").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = "funkar ej"

# B22 "påbörjat" -> light-blue Brödtext font, same as B14
$ws.Range("B14").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = "påbörjat"

# B23 / B24 "ej påbörjat" already carry the correct red style (B24 kept
# the original "AIC" row's formatting across the row-insert); copy it
# onto B23 as well so both rows match.
$ws.Range("B24").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "ej påbörjat"
$ws.Range("B24").Value = "ej påbörjat"

$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 5. Selection bookkeeping (matches the author's final cursor spot)
# -----------------------------------------------------------------
$ws.Range("C27").Select()

Write-Output "done"
